$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "17:52:22", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:52:30", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:52:40", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:52:51", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:53:01", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:53:12", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 87
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
        $cell.Style = "Normal"
    }
}
